$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value2
$text = $text.Replace("1000 Bs = 13.78 = 55688.77 pesos", "1000 Bs = 13.74 = 55387.77 pesos")
$text = $text.Replace("55688.77 pesos = 13.75 = 973.93 Bs", "55387.77 pesos = 13.62 = 959.21 Bs")
$cellA1.Value = $text

# --- Update "tasas" sheet rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 72.8
$ws2.Range("O10").Value = 4032.23
$ws2.Range("N12").Value = 4068
$ws2.Range("O12").Value = 70.45
